$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = '51.437.90'
$ws.Range("E2").Value = '  +1.12%  '

# Row 3
$ws.Range("D3").Value = '2.988.13'
$ws.Range("E3").Value = '  +1.83%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '380.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.88%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.97%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.544'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.35%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("E9").Value = '  +1.93%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.10%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.138'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.58%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0859'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.02%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.79%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.442.56'
$ws.Range("E14").Value = '  +1.51%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.10%  '

# Row 16
$ws.Range("D16").Value = '2.980.78'
$ws.Range("E16").Value = '  +1.67%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '11.29'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -7.83%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.988'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.68%  '

# Row 19
$ws.Range("D19").Value = '51.478.61'
$ws.Range("E19").Value = '  +1.28%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.02%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.71%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0962'
$ws.Range("E22").Value = '  +1.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.47%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.78'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.43%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.23'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.78%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.03%  '

# Row 27
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.75%  '

# Row 28
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.170'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.13%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.08%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.15%  '

# Row 31
$ws.Range("E31").Value = '  +1.72%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.95%  '

# Row 33
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '51.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.42%  '

# Row 34
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '34.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.97%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.02'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.53%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0443'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.31%  '

# Row 37
$ws.Range("E37").Value = '  +0.00%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.32%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.77%  '

# Row 40
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '129.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.89%  '

# Row 41
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.65%  '

# Row 42
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.69%  '

# Row 43
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.84'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.48%  '

# Row 44
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.83'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +12.38%  '

# Row 45
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.274'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.06%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.32%  '

# Row 47
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.17%  '

# Row 48
$ws.Range("E48").Value = '  +1.30%  '

# Row 49
$ws.Range("D49").Value = '2.043.40'
$ws.Range("E49").Value = '  +2.57%  '

# Row 50
$ws.Range("D50").Value = '3.282.36'
$ws.Range("E50").Value = '  +1.87%  '

# Row 51
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.532'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +18.23%  '
